$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from ParticipationSig" sheet to "Include #0" ---
$includeSheet = $wb.Worksheets.Item("Include from ParticipationSig")
$includeSheet.Name = "Include #0"

# --- 2. Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Make room for a new "Jurisdiction" row by copying the formatting of the
# last existing data row (row 14, which uses the plain bordered body style)
# onto the brand-new row 15 before the remaining rows are re-written.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Re-write rows 11-15 so that the previous rows 11-14 content (Description,
# Purpose, Copyright, Immutable) shifts down by one row, and the new
# "Jurisdiction" property row is inserted at row 11.
$ws.Range("A15").Value = "Immutable"
$ws.Range("B15").Value = "BooleanType[null]"

$ws.Range("A14").Value = "Copyright"
$ws.Range("B14").Value = ""

$ws.Range("A13").Value = "Purpose"
$ws.Range("B13").Value = ""

$ws.Range("A12").Value = "Description"
$ws.Range("B12").Value = "A set of codes specifying whether and how the participant has attested his participation through a signature - limited to values allowed in original CDA definition.`n`n**Note:** CDA Release One represented either an intended (``X``) or actual (``S``) authenticator. CDA Release Two only represents an actual authenticator, so has deprecated the value of ``X``."
# Recompute the row height instead of leaving an explicit custom height
# after writing the multi-line wrapped text.
$ws.Rows("12:12").AutoFit()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
